$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB: copy number/date format from BA1 to BB1 (header row, date style)
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("BB1").Value2 = 45986

# Rows 2-72: BB column duplicates the BA column value exactly for each row
$bbSame = [ordered]@{
  2 = 2.672233110627005
  3 = -0.2486584437591262
  4 = 2.48444986756347
  5 = 1.924160321525846
  6 = 1.36313906046999
  7 = -4.856442119446953
  8 = 1.575066711296429
  9 = 1.814138146527952
  10 = 0.4238544839342779
  11 = -0.1584284749490763
  12 = 0.4882046181590169
  13 = 1.53878081519332
  14 = 0.8692906535860487
  15 = 3.11595497587993
  16 = 0.3906876693375665
  17 = 1.221836833579857
  18 = 1.653194230429179
  19 = -1.198649426118308
  20 = 0.7741673943688596
  21 = 0.8244429407371285
  22 = -0.3344927361763723
  23 = 0.7295104243506501
  24 = 0.03704431742310987
  25 = 0.2220551683158618
  26 = 0.1661617031019347
  27 = 1.924586088852507
  28 = 0.5703441042800677
  29 = 0.5671096279522487
  30 = 0.6731221874372437
  31 = 0.8295169162459786
  32 = 0.363674885967896
  33 = 0.6363538952886927
  34 = 0.4303070273019074
  35 = 0.907245662456674
  36 = 0.7278878628511336
  37 = 0.5936080878907575
  38 = 0.3249989166702818
  39 = 2.043550613228959
  40 = 0.867287375484608
  41 = 0.6331942894404392
  42 = -0.05439614307451279
  43 = 0.1739459843577862
  44 = 1.255538557350434
  45 = 1.090483027535811
  46 = 1.226659036647675
  47 = -0.5
  48 = 0.9
  49 = 1
  50 = 1.1
  51 = -2.477834671711193
  52 = 0.2428240279789122
  53 = 0.8650544612728055
  54 = 0.05143518179183104
  55 = -2.588552528306963
  56 = 1.067142397791443
  57 = 0.9403228036019016
  58 = 1.348411706012428
  59 = -0.1190215178375666
  60 = 0.7996487817115536
  61 = 1.039074166251879
  62 = 0.6958556561364588
  63 = -1.636815679601384
  64 = -0.04357278727286484
  65 = 0.3050453114869214
  66 = 0.7386451510207621
  67 = 1.021259612058628
  68 = 0.5487159577757694
  69 = 0.6800497182067176
  70 = 0.1667761162031525
  71 = 0.5709980498538272
  72 = 1.176666004305858
}
foreach ($r in $bbSame.Keys) {
  $ws.Cells.Item($r, 54).Value2 = $bbSame[$r]
}

# Rows 73-82: BB column gets new distinct forecast values (EQUIPMENT eval update)
$bbNew = [ordered]@{
  73 = 0.8783323788356512
  74 = 0.375938928755442
  75 = 0.4360287204413526
  76 = 0.5352888176392425
  77 = 0.5231897612750926
  78 = 0.5035789880531426
  79 = 0.5060143287724851
  80 = 0.5098887247466805
  81 = 0.5093986862017086
  82 = 0.5086332614023801
}
foreach ($r in $bbNew.Keys) {
  $ws.Cells.Item($r, 54).Value2 = $bbNew[$r]
}

# New row 83: copy style of A82 (date format) down to A83, set new quarter date, and set BB83 value
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A83").Value2 = 46934
$ws.Cells.Item(83, 54).Value2 = 0.5087318365954954

$excel.CutCopyMode = 0

Write-Host "Done."